$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task list text (shared strings) in place, preserving order
$ws.Range("A1").Value = "create home page"
$ws.Range("A2").Value = "login stuff"
$ws.Range("A3").Value = "user dashboard"
$ws.Range("A4").Value = "admin dashboard"
$ws.Range("A5").Value = "event calander"

# Move the active selection from C4 to A5, as reflected in the sheet view
$ws.Range("A5").Select()
